# Update "想去人数" (want-to-go count) values in column F across all sheets
# to reflect the refreshed scrape data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2350
$ws.Range("F3").Value = 518
$ws.Range("F4").Value = 207
$ws.Range("F5").Value = 349
$ws.Range("F6").Value = 349
$ws.Range("F7").Value = 563
$ws.Range("F9").Value = 779
$ws.Range("F10").Value = 539
$ws.Range("F11").Value = 802
$ws.Range("F12").Value = 387
$ws.Range("F13").Value = 98
$ws.Range("F14").Value = 399
$ws.Range("F16").Value = 1026
$ws.Range("F17").Value = 20928
$ws.Range("F18").Value = 765
$ws.Range("F19").Value = 72
$ws.Range("F20").Value = 253
$ws.Range("F21").Value = 296
$ws.Range("F22").Value = 176
$ws.Range("F23").Value = 162
$ws.Range("F24").Value = 14
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 225
$ws.Range("F28").Value = 347

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 36
$ws.Range("F3").Value = 184
$ws.Range("F7").Value = 225
$ws.Range("F8").Value = 3427
$ws.Range("F10").Value = 98
$ws.Range("F16").Value = 3829

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 269
$ws.Range("F3").Value = 104
$ws.Range("F4").Value = 608
$ws.Range("F5").Value = 209

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 269
$ws.Range("F3").Value = 104
$ws.Range("F4").Value = 36
$ws.Range("F5").Value = 2350
$ws.Range("F6").Value = 608
$ws.Range("F7").Value = 518
$ws.Range("F8").Value = 207
$ws.Range("F9").Value = 349
$ws.Range("F10").Value = 349
$ws.Range("F11").Value = 563
$ws.Range("F12").Value = 184
$ws.Range("F17").Value = 209
$ws.Range("F18").Value = 779
$ws.Range("F19").Value = 539
$ws.Range("F20").Value = 802
$ws.Range("F21").Value = 387
$ws.Range("F22").Value = 98
$ws.Range("F23").Value = 399
$ws.Range("F25").Value = 1026
$ws.Range("F26").Value = 20928
$ws.Range("F27").Value = 225
$ws.Range("F28").Value = 3427
$ws.Range("F30").Value = 98
$ws.Range("F32").Value = 765
$ws.Range("F33").Value = 72
$ws.Range("F34").Value = 253
$ws.Range("F37").Value = 296
$ws.Range("F38").Value = 176
$ws.Range("F39").Value = 162
$ws.Range("F40").Value = 14
$ws.Range("F41").Value = 11
$ws.Range("F44").Value = 225
$ws.Range("F46").Value = 347
$ws.Range("F48").Value = 3829
